$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns with refreshed crypto data.
# For Price cells whose new value would otherwise be auto-parsed as a number
# by Excel (losing the original text formatting, e.g. trailing zeros or the
# thousands-grouped "63.272.86" style text), force the cell to Text first,
# then restore the cell style to Normal so no stray formatting is introduced.

$ws.Range('D2').Value = '63.272.86'
$ws.Range('E2').Value = '  +0.87%  '
$ws.Range('D3').Value = '2.646.56'
$ws.Range('E3').Value = '  +2.54%  '
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '594.60'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +2.31%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '142.99'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.95%  '
$ws.Range('E7').Value = '  +0.06%  '
$ws.Range('E8').Value = '  -0.82%  '
$ws.Range('D9').Value = '2.645.62'
$ws.Range('E9').Value = '  +2.53%  '
$ws.Range('E10').Value = '  -0.11%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '5.66'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +1.61%  '
$ws.Range('E12').Value = '  +0.80%  '
$ws.Range('E13').Value = '  +0.63%  '
$ws.Range('E14').Value = '  +1.55%  '
$ws.Range('D15').Value = '3.120.53'
$ws.Range('E15').Value = '  +2.52%  '
$ws.Range('D16').Value = '63.156.58'
$ws.Range('E16').Value = '  +0.78%  '
$ws.Range('E17').Value = '  -0.44%  '
$ws.Range('D18').Value = '2.641.22'
$ws.Range('E18').Value = '  +1.71%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '11.38'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +1.53%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '339.22'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +0.08%  '
$ws.Range('E21').Value = '  +0.26%  '
$ws.Range('E22').Value = '  +0.67%  '
$ws.Range('E23').Value = '  +0.01%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '66.87'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -0.79%  '
$ws.Range('E25').Value = '  +4.90%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.53'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -1.53%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.164'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -0.62%  '
$ws.Range('E28').Value = '  +0.42%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '8.42'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +2.10%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '7.79'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -2.77%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '527.83'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +15.64%  '
$ws.Range('E32').Value = '  +12.12%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.97'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +1.54%  '
$ws.Range('D34').Value = '0.0₃0804'
$ws.Range('E34').Value = '  -0.39%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '174.21'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -1.23%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '4.90'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +10.34%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.402'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +0.85%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '19.01'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +0.43%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.79'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +6.21%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '171.79'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +7.65%  '
$ws.Range('E42').Value = '  +0.02%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '40.21'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +0.83%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '3.73'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +0.55%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '21.97'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +4.53%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0558'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +4.29%  '
$ws.Range('E47').Value = '  -0.01%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0960'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -0.34%  '
$ws.Range('E49').Value = '  +1.49%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '18.50'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +2.20%  '
$ws.Range('E51').Value = '  -0.69%  '
